# Updates the cryptos worksheet Price (D) and Volume(1h) (E) columns
# to reflect the latest scrape, per the GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" values (column D) for the rows whose price changed.
$priceUpdates = @{
    2 = '37.795.87';
    3 = '2.080.36';
    5 = '233.18';
    7 = '58.57';
    9 = '0.393';
    10 = '0.0786';
    12 = '14.88';
    13 = '2.382.56';
    14 = '21.13';
    15 = '0.782';
    16 = '5.36';
    17 = '2.073.67';
    18 = '37.735.70';
    19 = '6.16';
    20 = '71.44';
    21 = '0.0₃0840';
    22 = '229.46';
    25 = '2.39';
    26 = '9.76';
    27 = '172.12';
    32 = '4.73';
    33 = '0.0633';
    34 = '4.68';
    37 = '3.41';
    38 = '0.999';
    39 = '5.42';
    40 = '0.0234';
    41 = '101.27';
    42 = '0.0974';
    44 = '16.91';
    45 = '1.448.52';
    46 = '1.16';
    47 = '1.07';
    48 = '4.10';
    49 = '7.35';
    51 = '2.269.50'
}

# New "Volume(1h)" values (column E) for every changed row.
$volumeUpdates = @{
    2 = '  -0.01%  ';
    3 = '  -0.26%  ';
    4 = '  -0.04%  ';
    5 = '  -0.41%  ';
    6 = '  -0.04%  ';
    7 = '  -0.34%  ';
    8 = '  -0.02%  ';
    9 = '  +0.31%  ';
    10 = '  -0.75%  ';
    11 = '  +3.22%  ';
    12 = '  +0.70%  ';
    13 = '  -0.50%  ';
    14 = '  -0.47%  ';
    15 = '  +1.66%  ';
    16 = '  +1.15%  ';
    17 = '  -0.54%  ';
    18 = '  +0.09%  ';
    19 = '  -1.52%  ';
    20 = '  -0.07%  ';
    21 = '  +1.24%  ';
    22 = '  +0.25%  ';
    24 = '  -0.59%  ';
    25 = '  +1.13%  ';
    26 = '  +7.89%  ';
    27 = '  +0.96%  ';
    28 = '  -1.58%  ';
    29 = '  -0.58%  ';
    30 = '  -1.45%  ';
    31 = '  +1.09%  ';
    32 = '  +0.45%  ';
    33 = '  +0.22%  ';
    34 = '  -0.96%  ';
    35 = '  -1.60%  ';
    36 = '  -0.80%  ';
    37 = '  -1.86%  ';
    38 = '  -0.19%  ';
    39 = '  +0.24%  ';
    40 = '  +8.74%  ';
    41 = '  +2.49%  ';
    42 = '  -0.62%  ';
    43 = '  -0.30%  ';
    44 = '  +4.90%  ';
    45 = '  -0.75%  ';
    46 = '  -1.57%  ';
    47 = '  -0.66%  ';
    48 = '  -4.42%  ';
    49 = '  -1.42%  ';
    50 = '  -1.80%  ';
    51 = '  -0.41%  '
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item([int]$row, 4)
    # Force text so values like "233.18" aren't auto-coerced to numbers
    # (matches the source data, which stores prices as plain text).
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item([int]$row, 5).Value = $volumeUpdates[$row]
}
